# Update countries & provincias Spain
#
# 1) Re-order "Hungria" in the country list so it sits right after "Armenia"
#    (i.e. before "Bulgaria"/"Taiwan"), and refresh its stats to the new
#    reported totals. The rows that used to hold Bulgaria/Taiwan/Hungria
#    keep their old Bulgaria/Taiwan numbers (they just "shift down" one
#    slot in country order), while the first of the three rows now carries
#    Hungria's fresh numbers.
# 2) Refresh a handful of other countries' case counters.
# 3) Bump the "last updated" timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder Armenia/Hungria/Bulgaria/Taiwan block (rows 68-71) ---------
# Row 68 = Armenia (untouched)
# Row 69 used to be Bulgaria -> becomes Hungria (with updated figures)
# Row 70 used to be Taiwan   -> becomes Bulgaria (old Bulgaria figures)
# Row 71 used to be Hungria  -> becomes Taiwan   (old Taiwan figures)

$ws.Range("A69").Value = "Hungria"
$ws.Range("B69").Value = 261
$ws.Range("C69").Value = 35
$ws.Range("D69").Value = 28
$ws.Range("E69").Value = 223
$ws.Range("F69").Value = 6
$ws.Range("G69").Value = 0
$ws.Range("H69").Value = 10

$ws.Range("A70").Value = "Bulgaria"
$ws.Range("B70").Value = 242
$ws.Range("C70").Value = 0
$ws.Range("D70").Value = 4
$ws.Range("E70").Value = 235
$ws.Range("F70").Value = 8
$ws.Range("G70").Value = 0
$ws.Range("H70").Value = 3

$ws.Range("A71").Value = "Taiwan"
$ws.Range("B71").Value = 235
$ws.Range("C71").Value = 0
$ws.Range("D71").Value = 29
$ws.Range("E71").Value = 204
$ws.Range("F71").Value = 0
$ws.Range("G71").Value = 0
$ws.Range("H71").Value = 2

# --- Refresh other countries' figures -----------------------------------

# Row 6: Estados Unidos
$ws.Range("B6").Value = 68489
$ws.Range("C6").Value = 278
$ws.Range("E6").Value = 67063
$ws.Range("F6").Value = 1455

# Row 18: Noruega
$ws.Range("B18").Value = 3100
$ws.Range("C18").Value = 16
$ws.Range("E18").Value = 3080

# Row 20: Australia
$ws.Range("B20").Value = 2766
$ws.Range("C20").Value = 90
$ws.Range("E20").Value = 2584

# Row 35: Tailandia
$ws.Range("B35").Value = 1045
$ws.Range("C35").Value = 111
$ws.Range("E35").Value = 971

# Row 44: India
$ws.Range("E44").Value = 609
$ws.Range("G44").Value = 1
$ws.Range("H44").Value = 13

# Row 100: Kazajistan
$ws.Range("D100").Value = 2
$ws.Range("E100").Value = 86

# Row 128: Macao
$ws.Range("B128").Value = 31
$ws.Range("E128").Value = 21

# --- Last updated timestamp ---------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 26 de Marzo de 2020 a las 06:42"
